$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "section header" rows ("situação do domicílio" and
# "grandes regiões") that no longer carry any data of their own; the
# rows below shift up to fill the gap, which realigns the already
# existing figures one row higher (e.g. what was the "urbana" figures
# row now directly follows "brasil", etc.).
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
